$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---- VENTAS POR GRUPO ----
$ws1.Range("L12").Value = 537.34
$ws1.Range("M20").Value = 1874.48
$ws1.Range("O20").Value = 323.48
$ws1.Range("L38").Value = 0
$ws1.Range("D55").Value = 457.92
$ws1.Range("M55").Value = -144.61
$ws1.Range("I58").Value = 356.35
$ws1.Range("M58").Value = 848.64
$ws1.Range("D78").Value = 475.2
$ws1.Range("L78").Value = 447.78
$ws1.Range("M78").Value = 685.4400000000001
$ws1.Range("M79").Value = 2422.91
$ws1.Range("M85").Value = 2056.32
$ws1.Range("N85").Value = 372.71
$ws1.Range("M92").Value = 1924.35
$ws1.Range("C100").Value = 518.4
$ws1.Range("L100").Value = 1089.41
$ws1.Range("M101").Value = 2683.55
$ws1.Range("C102").Value = 414.72
$ws1.Range("M102").Value = -636.72
$ws1.Range("M103").Value = 2762.6
$ws1.Range("D114").Value = 933.12
$ws1.Range("M114").Value = 1280.2
$ws1.Range("L117").Value = 855.36
$ws1.Range("M117").Value = 1382.34
$ws1.Range("D145").Value = 1373.76
$ws1.Range("M145").Value = 3259.2
$ws1.Range("O145").Value = 547.4299999999999
$ws1.Range("H153").Value = 811.8
$ws1.Range("O153").Value = 547.4299999999999
$ws1.Range("M218").Value = 2719.98
$ws1.Range("O244").Value = 831.1799999999999
$ws1.Range("L298").Value = 853.29
$ws1.Range("M298").Value = 5005.24
$ws1.Range("R299").Value = -10.44
$ws1.Range("M323").Value = -347.92
$ws1.Range("C327").Value = "6 de 325"
$ws1.Range("D327").Value = "23 de 325"
$ws1.Range("H327").Value = "12 de 325"
$ws1.Range("I327").Value = "16 de 325"
$ws1.Range("L327").Value = "29 de 325"
$ws1.Range("M327").Value = "60 de 325"
$ws1.Range("N327").Value = "4 de 325"
$ws1.Range("O327").Value = "6 de 325"

# ---- VENTA MENSUAL ----
$ws2.Range("F12").Value = 998.71
$ws2.Range("F20").Value = 2633.74
$ws2.Range("F38").Value = 0
$ws2.Range("F55").Value = 313.31
$ws2.Range("F58").Value = 1204.99
$ws2.Range("F78").Value = 2068.07
$ws2.Range("F79").Value = 2781.14
$ws2.Range("F85").Value = 2429.03
$ws2.Range("F92").Value = 4227.03
$ws2.Range("F100").Value = 1607.81
$ws2.Range("F101").Value = 2683.55
$ws2.Range("F102").Value = 1305.91
$ws2.Range("F103").Value = 2762.6
$ws2.Range("F114").Value = 2213.32
$ws2.Range("F117").Value = 2237.7
$ws2.Range("F149").Value = 5180.39
$ws2.Range("F157").Value = 5833.17
$ws2.Range("F222").Value = 2719.98
$ws2.Range("F248").Value = 831.1799999999999
$ws2.Range("F302").Value = 5858.53
$ws2.Range("F303").Value = -10.44
$ws2.Range("F327").Value = -347.92
$ws2.Range("F331").Value = 231951.69

# ---- CUMPLIMIENTO MENSUAL ----
$ws3.Range("D11").Value = 3761.38
$ws3.Range("E11").Value = -839.1554181472602
$ws3.Range("F11").Value = 1.287163219198992
$ws3.Range("D12").Value = 10084.97
$ws3.Range("E12").Value = 12348.7853751766
$ws3.Range("F12").Value = 0.449544440123441
$ws3.Range("D14").Value = 323.48
$ws3.Range("E14").Value = 1101.4862010375
$ws3.Range("F14").Value = 0.2270088930982912
$ws3.Range("D16").Value = 10969.54
$ws3.Range("E16").Value = 9417.937421713497
$ws3.Range("F16").Value = 0.5380528337614241
$ws3.Range("D20").Value = 518.9
$ws3.Range("E20").Value = 367.811016287574
$ws3.Range("F20").Value = 0.5851962933453764
$ws3.Range("D24").Value = 5268.7
$ws3.Range("E24").Value = 14304.3602492497
$ws3.Range("F24").Value = 0.2691812078901646
$ws3.Range("D25").Value = 19863.63
$ws3.Range("E25").Value = 27270.6331579098
$ws3.Range("F25").Value = 0.4214265519215314
$ws3.Range("D26").Value = 993.89
$ws3.Range("E26").Value = 116.5466512034101
$ws3.Range("F26").Value = 0.8950443043490096
$ws3.Range("D28").Value = 933.12
$ws3.Range("E28").Value = 5264.46402943659
$ws3.Range("F28").Value = 0.1505618956625632
$ws3.Range("D29").Value = 2697.4
$ws3.Range("E29").Value = 14971.7470988183
$ws3.Range("F29").Value = 0.1526615849035747
$ws3.Range("D37").Value = 9691.5
$ws3.Range("E37").Value = 8139.9143984654
$ws3.Range("F37").Value = 0.5435070815713904
$ws3.Range("D38").Value = 23143.39
$ws3.Range("E38").Value = 38720.3303947566
$ws3.Range("F38").Value = 0.3741027835429304
$ws3.Range("D42").Value = 1848.96
$ws3.Range("E42").Value = 3655.65890386263
$ws3.Range("F42").Value = 0.3358924627284501
$ws3.Range("D45").Value = 2359.65
$ws3.Range("E45").Value = 547.9336814602598
$ws3.Range("F45").Value = 0.8115501593456894
$ws3.Range("D51").Value = 38098.83
$ws3.Range("E51").Value = -1275.186907882904
$ws3.Range("F51").Value = 1.034629569504922
$ws3.Range("D53").Value = 1094.86
$ws3.Range("E53").Value = -178.761404787216
$ws3.Range("F53").Value = 1.195133368527538
$ws3.Range("D63").Value = 238.82
$ws3.Range("E63").Value = 3267.84949822329
$ws3.Range("F63").Value = 0.0681045077447424
$ws3.Range("D64").Value = 18867.86
$ws3.Range("E64").Value = 13536.94
$ws3.Range("F64").Value = 0.5822550980101714
$ws3.Range("D79").Value = 831.1799999999999
$ws3.Range("E79").Value = -163.496851612446
$ws3.Range("F79").Value = 1.244871915679299
$ws3.Range("D80").Value = 8761.92
$ws3.Range("E80").Value = -8761.92
$ws3.Range("D81").Value = 11551.89
$ws3.Range("E81").Value = 8448.110000000001
$ws3.Range("F81").Value = 0.5775945
$ws3.Range("D92").Value = 19447.67
$ws3.Range("E92").Value = 23652.4154117774
$ws3.Range("F92").Value = 0.4512211475730808
$ws3.Range("D95").Value = -10.44
$ws3.Range("E95").Value = 10.44
$ws3.Range("D97").Value = 242145.0899999999
$ws3.Range("E97").Value = 253370.5006021116
$ws3.Range("F97").Value = 0.4886729995836545
